$p = $ppt.ActivePresentation
$s = $p.Slides.Item(10)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

# Paragraph with "$ python extract_email.py" is the anchor we extend from,
# so the new runs inherit its exact bold/lang/dirty run formatting.
$anchor = $tr.Paragraphs(5, 1)
$anchor.InsertAfter("`rDifferent Email Formats To Cover:`rEugene dot agichten at emory dot edu`rRohini [@] buffalo [DOT] edu")

# Bump the indent level of the two new "example address" bullets to lvl="2"
# (COM IndentLevel is 1-based, so XML lvl="2" == IndentLevel 3).
$tr.Paragraphs(7, 1).IndentLevel = 3
$tr.Paragraphs(8, 1).IndentLevel = 3
